$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.1636683333333333
$ws.Range("H2").Value = 0.491005
$ws.Range("I2").Value = 0.008639493057305454
$ws.Range("J2").Value = 0.008639493057305455
$ws.Range("M2").Value = 24.75542533333333
$ws.Range("N2").Value = 74.26627599999999
$ws.Range("O2").Value = 0.7762421087066456
$ws.Range("P2").Value = 0.7762421087066456
$ws.Range("Q2").Value = 4.051679205264444
$ws.Range("R2").Value = 36.46511284738
$ws.Range("S2").Value = 0.00670633830895921
$ws.Range("T2").Value = 0.006706338308959211

$ws.Range("G3").Value = 0.1636683333333333
$ws.Range("H3").Value = 0.491005
$ws.Range("I3").Value = 0.008639493057305454
$ws.Range("J3").Value = 0.008639493057305455
$ws.Range("M3").Value = 3.818542
$ws.Range("O3").Value = 0.1197358984688377
$ws.Range("P3").Value = 0.1197358984688377
$ws.Range("Q3").Value = 0.6249744049033333
$ws.Range("R3").Value = 5.624769644130001
$ws.Range("S3").Value = 0.001034457463531754
$ws.Range("T3").Value = 0.001034457463531754

$ws.Range("G4").Value = 0.1636683333333333
$ws.Range("H4").Value = 0.491005
$ws.Range("I4").Value = 0.008639493057305454
$ws.Range("J4").Value = 0.008639493057305455
$ws.Range("M4").Value = 3.317404
$ws.Range("N4").Value = 9.952211999999999
$ws.Range("O4").Value = 0.1040219928245168
$ws.Range("P4").Value = 0.1040219928245168
$ws.Range("Q4").Value = 0.5429539836733333
$ws.Range("R4").Value = 4.88658585306
$ws.Range("S4").Value = 0.0008986972848144903
$ws.Range("T4").Value = 0.0008986972848144906

$ws.Range("I5").Value = 0.808839719627903
$ws.Range("J5").Value = 0.8088397196279031
$ws.Range("M5").Value = 24.75542533333333
$ws.Range("N5").Value = 74.26627599999999
$ws.Range("O5").Value = 0.7762421087066456
$ws.Range("P5").Value = 0.7762421087066456
$ws.Range("Q5").Value = 379.3230749386586
$ws.Range("R5").Value = 3413.907674447928
$ws.Range("S5").Value = 0.6278554495696554
$ws.Range("T5").Value = 0.6278554495696556

$ws.Range("I6").Value = 0.808839719627903
$ws.Range("J6").Value = 0.8088397196279031
$ws.Range("M6").Value = 3.818542
$ws.Range("O6").Value = 0.1197358984688377
$ws.Range("P6").Value = 0.1197358984688377
$ws.Range("S6").Value = 0.09684715054692973
$ws.Range("T6").Value = 0.09684715054692974

$ws.Range("I7").Value = 0.808839719627903
$ws.Range("J7").Value = 0.8088397196279031
$ws.Range("M7").Value = 3.317404
$ws.Range("N7").Value = 9.952211999999999
$ws.Range("O7").Value = 0.1040219928245168
$ws.Range("P7").Value = 0.1040219928245168
$ws.Range("Q7").Value = 50.832004263704
$ws.Range("R7").Value = 457.488038373336
$ws.Range("S7").Value = 0.08413711951131789
$ws.Range("T7").Value = 0.0841371195113179

$ws.Range("G8").Value = 3.457711333333334
$ws.Range("H8").Value = 10.373134
$ws.Range("I8").Value = 0.1825207873147914
$ws.Range("J8").Value = 0.1825207873147914
$ws.Range("M8").Value = 24.75542533333333
$ws.Range("N8").Value = 74.26627599999999
$ws.Range("O8").Value = 0.7762421087066456
$ws.Range("P8").Value = 0.7762421087066456
$ws.Range("Q8").Value = 85.59711473655378
$ws.Range("R8").Value = 770.374032628984
$ws.Range("S8").Value = 0.1416803208280309
$ws.Range("T8").Value = 0.1416803208280309

$ws.Range("G9").Value = 3.457711333333334
$ws.Range("H9").Value = 10.373134
$ws.Range("I9").Value = 0.1825207873147914
$ws.Range("J9").Value = 0.1825207873147914
$ws.Range("M9").Value = 3.818542
$ws.Range("O9").Value = 0.1197358984688377
$ws.Range("P9").Value = 0.1197358984688377
$ws.Range("Q9").Value = 13.20341595020933
$ws.Range("R9").Value = 118.830743551884
$ws.Range("S9").Value = 0.02185429045837618
$ws.Range("T9").Value = 0.02185429045837618

$ws.Range("G10").Value = 3.457711333333334
$ws.Range("H10").Value = 10.373134
$ws.Range("I10").Value = 0.1825207873147914
$ws.Range("J10").Value = 0.1825207873147914
$ws.Range("M10").Value = 3.317404
$ws.Range("N10").Value = 9.952211999999999
$ws.Range("O10").Value = 0.1040219928245168
$ws.Range("P10").Value = 0.1040219928245168
$ws.Range("Q10").Value = 11.47062540804533
$ws.Range("R10").Value = 103.235628672408
$ws.Range("S10").Value = 0.01898617602838439
$ws.Range("T10").Value = 0.01898617602838439
